$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("B2").Value = 0.04172184405617529
$ws.Range("C2").Value = 0.002658071450198252
$ws.Range("D2").Value = 189.6080260415259
$ws.Range("E2").Value = 13.86384647080068
$ws.Range("G2").Value = 203.516252427833
